$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("C2").Value = 40
$ws.Range("D2").Value = 3

# Add new rows 3-6
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "run-01"
$ws.Range("C3").Value = 46
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = "--"
$ws.Range("J3").Value = "test"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "run-01"
$ws.Range("C4").Value = 60
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = "--"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "run-01"
$ws.Range("C5").Value = 167
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = "--"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "run-01"
$ws.Range("C6").Value = 240
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = "--"
